# Revert "Update to Dumbarton Rail vehicle type and capacity"
# - transitPrefixToVehicle!C19: "DBRail 4 car" -> "Dumbarton Rail"
# - transitLineToVehicle: remove the "134_DBRail" data row (was row 498)
# - transitVehicleToCapacity: remove the "DBRail 4 car" vehicle row (was row 21)

$wb = $excel.ActiveWorkbook

# --- transitVehicleToCapacity: delete the DBRail 4 car row ---
$ws3 = $wb.Worksheets.Item("transitVehicleToCapacity")
$ws3.Rows.Item(21).Delete()

# --- transitLineToVehicle: delete the 134_DBRail line row ---
$ws2 = $wb.Worksheets.Item("transitLineToVehicle")
$ws2.Rows.Item(498).Delete()

# --- transitPrefixToVehicle: restore the old vehicle type value ---
$ws1 = $wb.Worksheets.Item("transitPrefixToVehicle")
$ws1.Range("C19").Value = "Dumbarton Rail"

# --- restore on-screen selections to match the reverted state ---
$ws1.Activate()
$ws1.Range("F21").Select()

$ws3.Activate()
$ws3.Range("A55").Select()

$ws2.Activate()
$ws2.Range("E21").Select()
